$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Well Authorization Number" column of data between the existing
# "Unique ID" and "Amount" columns (columns D/E) of the payment details table,
# without disturbing the overall table dimensions (still A1:E15).
#
# For each affected row we first copy the current column D content into column
# E (this is where "Amount" / its placeholders used to live), then overwrite
# column D with the new Well Authorization Number header / placeholders.

# Header row (row 11): D11 "Amount" -> E11, D11 becomes "Well Authorization Number"
$ws.Range("E11").Value2 = $ws.Range("D11").Value2
$ws.Range("D11").Value2 = "Well Authorization Number"

# Data row 12: D12 "{d.payment_details[i].amount}" -> E12, D12 becomes the WAN placeholder
$ws.Range("E12").Value2 = $ws.Range("D12").Value2
$ws.Range("D12").Value2 = "{d.payment_details[i].well_authorization_number}"

# Data row 13: D13 "{d.payment_details[i+1].amount}" -> E13, D13 becomes the WAN placeholder
$ws.Range("E13").Value2 = $ws.Range("D13").Value2
$ws.Range("D13").Value2 = "{d.payment_details[i+1].well_authorization_number}"

# Totals row 14: D14 "{d.total_payment}" -> E14, D14 becomes a blank spacer (" ")
$ws.Range("E14").Value2 = $ws.Range("D14").Value2
$ws.Range("D14").Value2 = " "

# Update the saved selection to match what was active when the workbook was saved.
$ws.Range("D18").Select()
